$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 'HİLMİ MÜFTÜOĞLU, SÜLEYMAN YILMAZ, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E3").Value = 'HİLMİ MÜFTÜOĞLU, SÜLEYMAN YILMAZ, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E4").Value = 'HİLMİ MÜFTÜOĞLU, SÜLEYMAN YILMAZ, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E5").Value = 'HİLMİ MÜFTÜOĞLU, SÜLEYMAN YILMAZ, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E6").Value = 'HİLMİ MÜFTÜOĞLU, SÜLEYMAN YILMAZ, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E7").Value = 'HİLMİ MÜFTÜOĞLU, SÜLEYMAN YILMAZ, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E8").Value = 'HİLMİ MÜFTÜOĞLU, SÜLEYMAN YILMAZ, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E9").Value = 'HİLMİ MÜFTÜOĞLU, SÜLEYMAN YILMAZ, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E10").Value = 'HİLMİ MÜFTÜOĞLU, SÜLEYMAN YILMAZ, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E33").Value = 'HİLMİ MÜFTÜOĞLU, ÖMER ÇEVİK, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E34").Value = 'HİLMİ MÜFTÜOĞLU, ÖMER ÇEVİK, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E35").Value = 'HİLMİ MÜFTÜOĞLU, ÖMER ÇEVİK, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E36").Value = 'HİLMİ MÜFTÜOĞLU, ÖMER ÇEVİK, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E37").Value = 'HİLMİ MÜFTÜOĞLU, ÖMER ÇEVİK, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E38").Value = 'HİLMİ MÜFTÜOĞLU, ÖMER ÇEVİK, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E39").Value = 'HİLMİ MÜFTÜOĞLU, ÖMER ÇEVİK, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E40").Value = 'HİLMİ MÜFTÜOĞLU, AHMET TÜRKOĞLU, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E41").Value = 'HİLMİ MÜFTÜOĞLU, AHMET TÜRKOĞLU, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E42").Value = 'HİLMİ MÜFTÜOĞLU, AHMET TÜRKOĞLU, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E43").Value = 'HİLMİ MÜFTÜOĞLU, AHMET TÜRKOĞLU, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E44").Value = 'HİLMİ MÜFTÜOĞLU, AHMET TÜRKOĞLU, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E45").Value = 'HİLMİ MÜFTÜOĞLU, AHMET TÜRKOĞLU, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E46").Value = 'HİLMİ MÜFTÜOĞLU, AHMET TÜRKOĞLU, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E47").Value = 'HİLMİ MÜFTÜOĞLU, AHMET TÜRKOĞLU, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E48").Value = 'HİLMİ MÜFTÜOĞLU, AYHAN KARADAYI, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E49").Value = 'HİLMİ MÜFTÜOĞLU, AYHAN KARADAYI, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E50").Value = 'HİLMİ MÜFTÜOĞLU, AYHAN KARADAYI, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E51").Value = 'HİLMİ MÜFTÜOĞLU, AYHAN KARADAYI, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E52").Value = 'HİLMİ MÜFTÜOĞLU, AYHAN KARADAYI, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E53").Value = 'HİLMİ MÜFTÜOĞLU, AYHAN KARADAYI, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E54").Value = 'HİLMİ MÜFTÜOĞLU, AYHAN KARADAYI, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E55").Value = 'HİLMİ MÜFTÜOĞLU, AYHAN KARADAYI, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E56").Value = 'HİLMİ MÜFTÜOĞLU, AYHAN KARADAYI, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E57").Value = 'HİLMİ MÜFTÜOĞLU, AYHAN KARADAYI, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E58").Value = 'HİLMİ MÜFTÜOĞLU, AYHAN KARADAYI, MEHMET AKDENİZ, MUSTAFA BICAK'
$ws.Range("E59").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E60").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E61").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E62").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E63").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E64").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E65").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E66").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E67").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E68").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E69").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E70").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E71").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E72").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E73").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E74").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E75").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E76").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E77").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E78").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E79").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E80").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E81").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E82").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E83").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E84").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E85").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E86").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E87").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E88").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E89").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E90").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E91").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E92").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E93").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E94").Value = 'KEMAL KORKMAZ, ENGİN UĞURLU, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E95").Value = 'CİHAN KARA, MUSTAFA GÜRBÜZ, ENDER NUSRET ÖNAL GÜLSOY, İSMAİL AKLAN'
$ws.Range("E96").Value = 'MEHMET NEJAT AY, HASAN İN, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E97").Value = 'MEHMET NEJAT AY, HASAN İN, ENDER NUSRET ÖNAL GÜLSOY, İSMAİL AKLAN'
$ws.Range("E98").Value = 'MEHTAP AKDOĞAN, CANER OKAY, MEHMET AKDENİZ, MUSTAFA BICAK, MÜMİN AKDOĞAN'
$ws.Range("E99").Value = 'MEHMET NEJAT AY, MUSTAFA GÜRBÜZ, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E100").Value = 'MEHTAP AKDOĞAN, CANER OKAY, MEHMET AKDENİZ, MUSTAFA BICAK, MÜMİN AKDOĞAN'
$ws.Range("E101").Value = 'MEHMET NEJAT AY, ŞEKİP KORKMAZ, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E102").Value = 'TEVFİK YILDIZ, HASAN İN, ENDER NUSRET ÖNAL GÜLSOY, İSMAİL AKLAN'
$ws.Range("E103").Value = 'MEHMET NEJAT AY, HASAN İN, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E104").Value = 'TEVFİK YILDIZ, HASAN İN, ENDER NUSRET ÖNAL GÜLSOY, İSMAİL AKLAN'
$ws.Range("E105").Value = 'TAHA GÜRKAN, MAHMUT MELEMEN, İSMAİL AKLAN, TEMUR ARSLAN'
$ws.Range("E106").Value = 'TAHA GÜRKAN, MAHMUT MELEMEN, İSMAİL AKLAN, TEMUR ARSLAN'
$ws.Range("E107").Value = 'CİHAN KARA, MAHMUT MELEMEN, ENDER NUSRET ÖNAL GÜLSOY, İSMAİL AKLAN'
$ws.Range("E108").Value = 'TAHA GÜRKAN, BARIŞ YAYLAGÜL, İSMAİL AKLAN, TEMUR ARSLAN'
$ws.Range("E109").Value = 'TAHA GÜRKAN, BARIŞ YAYLAGÜL, İSMAİL AKLAN, TEMUR ARSLAN'
$ws.Range("E110").Value = 'CİHAN KARA, MAHMUT MELEMEN, İSMAİL AKLAN, TEMUR ARSLAN'
$ws.Range("E111").Value = 'CİHAN KARA, MAHMUT MELEMEN, ENDER NUSRET ÖNAL GÜLSOY, İSMAİL AKLAN'
$ws.Range("E112").Value = 'TAHA GÜRKAN, BARIŞ YAYLAGÜL, İSMAİL AKLAN, TEMUR ARSLAN'
$ws.Range("E113").Value = 'TAHA GÜRKAN, BARIŞ YAYLAGÜL, İSMAİL AKLAN, TEMUR ARSLAN'
$ws.Range("E114").Value = 'MEHMET AKGÜN KOLUKIRIK, BARIŞ YAYLAGÜL, ENDER NUSRET ÖNAL GÜLSOY, İSMAİL AKLAN'
$ws.Range("E115").Value = 'TAHA GÜRKAN, MAHMUT MELEMEN, ENDER NUSRET ÖNAL GÜLSOY, İSMAİL AKLAN'
$ws.Range("E116").Value = 'MEHMET NEJAT AY, ŞEKİP KORKMAZ, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E117").Value = 'LOKMAN ALKAN, HASAN İN, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E118").Value = 'MEHMET NEJAT AY, ŞEKİP KORKMAZ, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E119").Value = 'CİHAN KARA, MUSTAFA GÜRBÜZ, ENDER NUSRET ÖNAL GÜLSOY, İSMAİL AKLAN'
$ws.Range("E120").Value = 'CİHAN KARA, ŞEKİP KORKMAZ, ENDER NUSRET ÖNAL GÜLSOY, İSMAİL AKLAN'
$ws.Range("E121").Value = 'MEHMET NEJAT AY, ŞEKİP KORKMAZ, ALİ BAŞKURT, İSMAİL AKLAN'
$ws.Range("E122").Value = 'MEHTAP AKDOĞAN, CANER OKAY, MEHMET AKDENİZ, MUSTAFA BICAK, MÜMİN AKDOĞAN'
